$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New header labels N1:S1 (training-set metric columns) ---
$newHeaders = @("MAE_tr", "MSE_tr", "RMSE_tr", "MAPE_tr", "ME_tr", "r2_tr")
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $col = 14 + $i   # N=14 .. S=19
    $ws.Cells.Item(1, $col).Value = $newHeaders[$i]
}

# Give the new headers the same bold/border/center style as the existing header row
$ws.Range("B1").Copy() | Out-Null
$ws.Range("N1:S1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# --- 2. Populate the new training-metric values for N2:S21 ---
$trainingData = @{
    2 = @(0.84073991661615; 1.206581138142671; 1.098444872600656; 47.80761582091964; 0.687594334313075; 0.7056280266967966)
    3 = @(0.7482752781868885; 0.9600988207505629; 0.9798463250686624; 43.14843560774903; 0.7514130614676118; 0.7516853382546603)
    4 = @(0.8219168437550652; 1.056373837791607; 1.027800485401523; 55.34996481796306; 0.7264857193793495; 0.7275187311624185)
    5 = @(0.8208388397635635; 1.082153572698884; 1.04026610667602; 53.5134702594744; 0.7198108800417001; 0.7198108800417131)
    6 = @(0.7799326383928791; 1.011933886138172; 1.005949246303298; 45.20183560991457; 0.7379920261170437; 0.7385794057947042)
    7 = @(0.515300722706108; 0.5230408314228558; 0.7232156189013452; 29.37309512338503; 0.8645752747522407; 0.8649174378973892)
    8 = @(0.5153019619756993; 0.523040831903816; 0.7232156192338604; 29.37339932756181; 0.8645752746277114; 0.8649175911434628)
    9 = @(0.4762441332043752; 0.4151614097185087; 0.6443302023950986; 33.33105519350939; 0.8925072069580978; 0.8925072069580979)
    10 = @(0.4153102562113989; 0.3558132028723484; 0.5965007987189526; 24.31615475203516; 0.9078735304327382; 0.9082158213071775)
    11 = @(0.357961232178416; 0.3292947756775786; 0.5738421173786207; 14.64294253496036; 0.914739630555524; 0.9180666752008412)
    12 = @(0.4337408261349388; 0.360596007991851; 0.6004964679262077; 27.61843420634286; 0.906635175737828; 0.9081458307018819)
    13 = @(0.459471490957871; 0.4249045799458637; 0.651847052571279; 24.02478892724872; 0.889984524077887; 0.890692786485581)
    14 = @(0.4670767667725075; 0.4229570038425562; 0.6503514464061384; 31.21938660864355; 0.8904887867335807; 0.8904887867336009)
    15 = @(0.4100463198398805; 0.3553973331531641; 0.5961521057189718; 21.18608094630648; 0.9079812066199036; 0.9079812066199099)
    16 = @(0.3574293354302711; 0.2536811070228989; 0.503667655327299; 25.30097764655129; 0.934317376091525; 0.9358909893881214)
    17 = @(0.761172496056864; 1.052843284620455; 1.026081519481008; 35.51931368350997; 0.7273998434103073; 0.7741805365716017)
    18 = @(0.7657950311563995; 1.066581049349283; 1.032754108851319; 34.69047875668944; 0.7238428878111454; 0.7751872251285193)
    19 = @(0.6299675668428945; 0.6593448709953675; 0.8120005363270196; 36.98472989593338; 0.8292836952037526; 0.8293849910394926)
    20 = @(0.4007662114371423; 0.33129592568787; 0.5755831179663542; 24.123925169314; 0.9142214966469607; 0.9142221926309912)
    21 = @(0.4138709543993681; 0.3299913284816157; 0.57444871701625; 25.73323765274334; 0.9145592804443277; 0.9148872467197487)
}
foreach ($r in $trainingData.Keys) {
    $rowVals = $trainingData[$r]
    for ($i = 0; $i -lt $rowVals.Length; $i++) {
        $ws.Cells.Item($r, 14 + $i).Value = $rowVals[$i]
    }
}

# --- 3. Correct the calibration metrics for row 10 (H10:M10) ---
$h10 = @(0.423513149854946; 0.3772584284019147; 0.6142136667332588; 24.04072077650301; 0.9008550587071262; 0.9016711456998099)
for ($i = 0; $i -lt $h10.Length; $i++) {
    $ws.Cells.Item(10, 8 + $i).Value = $h10[$i]
}
